# Updated symbol list on Sat Dec 24 15:06:36 UTC 2022 with GitHub Actions
#
# Applies the per-cell value updates produced by the refreshed crypto
# data feed: the "Hora" (Hour) column moves from 14 -> 15 for every
# data row, "Price" values are refreshed, and rows 9-17 roll the
# Coin/Link/Volume columns up by one position (with a new entry
# appearing at the bottom of that block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.48"
$ws.Range("G2").Value = "'15"
$ws.Range("D3").Value = "'21.89"
$ws.Range("G3").Value = "'15"
$ws.Range("D4").Value = "'5.391"
$ws.Range("G4").Value = "'15"
$ws.Range("D5").Value = "'0.05995"
$ws.Range("G5").Value = "'15"
$ws.Range("D6").Value = "'3.392"
$ws.Range("G6").Value = "'15"
$ws.Range("D7").Value = "'0.8150"
$ws.Range("G7").Value = "'15"
$ws.Range("D8").Value = "'0.9536"
$ws.Range("G8").Value = "'15"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1427"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("G9").Value = "'15"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07423"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("G10").Value = "'15"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03258"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").Value = "'15"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03057"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("G12").Value = "'15"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09417"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("G13").Value = "'15"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'4.003"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("G14").Value = "'15"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001592"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("G15").Value = "'15"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04813"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("G16").Value = "'15"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005912"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("G17").Value = "'15"
$ws.Range("D18").Value = "'0.005638"
$ws.Range("G18").Value = "'15"
$ws.Range("D19").Value = "'0.004151"
$ws.Range("G19").Value = "'15"
$ws.Range("D20").Value = "'0.0009923"
$ws.Range("G20").Value = "'15"
$ws.Range("G21").Value = "'15"
$ws.Range("D22").Value = "'3.671"
$ws.Range("G22").Value = "'15"
$ws.Range("D23").Value = "'6.427"
$ws.Range("G23").Value = "'15"
$ws.Range("D24").Value = "'2.189"
$ws.Range("G24").Value = "'15"
$ws.Range("G25").Value = "'15"
$ws.Range("G26").Value = "'15"
$ws.Range("G27").Value = "'15"
$ws.Range("G28").Value = "'15"
$ws.Range("G29").Value = "'15"
$ws.Range("G30").Value = "'15"
$ws.Range("G31").Value = "'15"
$ws.Range("G32").Value = "'15"
$ws.Range("G33").Value = "'15"
$ws.Range("G34").Value = "'15"
$ws.Range("G35").Value = "'15"
$ws.Range("G36").Value = "'15"
$ws.Range("G37").Value = "'15"
$ws.Range("G38").Value = "'15"
$ws.Range("G39").Value = "'15"
$ws.Range("D40").Value = "'0.03998"
$ws.Range("G40").Value = "'15"
$ws.Range("D41").Value = "'0.003048"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'15"
$ws.Range("G42").Value = "'15"
$ws.Range("D43").Value = "'0.002721"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("G43").Value = "'15"
$ws.Range("D44").Value = "'0.005769"
$ws.Range("G44").Value = "'15"
$ws.Range("D45").Value = "'0.00005135"
$ws.Range("G45").Value = "'15"
$ws.Range("G46").Value = "'15"
$ws.Range("D47").Value = "'0.8602"
$ws.Range("G47").Value = "'15"
$ws.Range("D48").Value = "'0.005598"
$ws.Range("G48").Value = "'15"
$ws.Range("G49").Value = "'15"
$ws.Range("G50").Value = "'15"
$ws.Range("G51").Value = "'15"
